# Edit script: update column G ("K") values for rows 2..74 in Sheet1
# per commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$kValues = @{
    2 = 0
    3 = 4
    4 = 1
    5 = 1
    6 = 1
    7 = 0
    8 = 1
    9 = 2
    10 = 1
    11 = 2
    12 = 2
    13 = 0
    14 = 3
    15 = 2
    16 = 1
    17 = 3
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    22 = 2
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    28 = 0
    29 = 2
    30 = 0
    31 = 1
    32 = 1
    33 = 2
    34 = 0
    35 = 0
    36 = 1
    37 = 2
    38 = 2
    39 = 2
    40 = 1
    41 = 0
    42 = 0
    43 = 0
    44 = 3
    45 = 0
    46 = 0
    47 = 0
    48 = 2
    49 = 2
    50 = 2
    51 = 2
    52 = 2
    53 = 0
    54 = 1
    55 = 0
    56 = 1
    57 = 0
    58 = 1
    59 = 1
    60 = 0
    61 = 1
    62 = 2
    63 = 0
    64 = 1
    65 = 0
    66 = 0
    67 = 2
    68 = 2
    69 = 1
    70 = 1
    71 = 1
    72 = 1
    73 = 2
    74 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
